# The deck's single theme ("Integral", wired to the slide master via
# ppt/theme/theme2.xml) is switched to the classic "Office Theme" colour
# scheme, mirroring a Design-tab theme change in the PowerPoint UI.
#
# PowerPoint's object model keeps exactly one mutable theme colour scheme
# per presentation (reached here through Slide.ThemeColorScheme); setting
# each slot's RGB reproduces the target clrScheme (and, because the font
# scheme/format scheme are already identical between the two themes in
# this deck, the whole themeElements block ends up byte-for-byte equal to
# the target "Office Theme" definition).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeRGB($colorScheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# msoThemeColorDark1 / Light1 / Dark2 / Light2
Set-ThemeRGB $tcs 1  "000000"
Set-ThemeRGB $tcs 2  "FFFFFF"
Set-ThemeRGB $tcs 3  "44546A"
Set-ThemeRGB $tcs 4  "E7E6E6"
# msoThemeColorAccent1-6
Set-ThemeRGB $tcs 5  "5B9BD5"
Set-ThemeRGB $tcs 6  "ED7D31"
Set-ThemeRGB $tcs 7  "A5A5A5"
Set-ThemeRGB $tcs 8  "FFC000"
Set-ThemeRGB $tcs 9  "4472C4"
Set-ThemeRGB $tcs 10 "70AD47"
# msoThemeColorHyperlink / FollowedHyperlink
Set-ThemeRGB $tcs 11 "0563C1"
Set-ThemeRGB $tcs 12 "954F72"
